# "Added last minute updates"
# - Paragraph 1 ("**ID__AFFARS_5322_topic_9__ID** "): the placeholder id is
#   updated to 406_13, the two runs are merged into a single run (dropping
#   the stray trailing space), a paragraph border is added, and the left
#   indent grows from 120 -> 225 twips (matching paragraph 3's formatting).

$d = $word.ActiveDocument

# 1. Replace the ID placeholder text. Matching the whole "<id>" + trailing
#    space also merges the two identically-formatted runs into one run,
#    exactly as in the target (no more stray trailing-space run).
$d.Content.Find.Execute("**ID__AFFARS_5322_topic_9__ID** ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_5322_406_13__ID**", 2)

# 2. Update paragraph formatting on paragraph 1: left indent 120 -> 225
#    twips (225 twips = 11.25 pt), and add a paragraph border with 5-twip
#    spacing on all four sides.
$p = $d.Paragraphs.Item(1)
$pf = $p.Range.ParagraphFormat
$pf.LeftIndent = 11.25

$borders = $pf.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
